# Daily refresh of the cryptos list: update Price (D) and Volume(1h) (E) columns
# for each coin row. Numeric-looking price strings are written with a leading
# apostrophe and then the cell style is reset to "Normal" so the value is stored
# as text (matching the original inlineStr formatting) without adding a visible
# quote-prefix indicator or changing the cell's style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.324.60"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.487.40"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D5").Value = "'596.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'177.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").Value = "'7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'0.426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "4.091.84"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "'31.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.05%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "67.292.49"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "3.483.97"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "'6.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "'14.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'388.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "'7.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'73.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'0.536"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "'5.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'10.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "'0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'23.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "'7.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "'1.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'163.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").Value = "'0.870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").Value = "'6.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").Value = "'27.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'4.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").Value = "'26.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "2.819.91"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'0.0722"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.14%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").Value = "'42.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'342.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "'33.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'6.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.97%  "
